# SinghRamesh_2023 metadata file - unit fixes for leaf N/P area variables
#
# 1. Row 41 (was "Leaf_C_area"): repurposed as "Leaf_P_Area_kg_m-2"
# 2. Row 42 (was "Leaf_N_Area"): repurposed as "Leaf_N_Area_kg_m-2", with a
#    description that records how the value is calculated
# 3. A brand-new row is inserted (new row 51) for "%P" / "Total leaf
#    phosphorous (%)", pushing the root isotope rows down by one
# 4. The numeric "pi" constant in B34 is given a left-aligned style
# 5. Column B is widened to fit the long description text
# 6. Final selection left on B55, matching where the author's cursor ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: fix the mislabeled leaf nutrient-per-area rows in place ---
$ws.Range("A41").Value = "Leaf_P_Area_kg_m-2"
$ws.Range("B41").Value = "Leaf phosphorous content per unit leaf area (calculated as : ((mass of Phosphorous/SLA)*1000) in kg m^2)"

$ws.Range("A42").Value = "Leaf_N_Area_kg_m-2"
$ws.Range("B42").Value = "Leaf nitrogen content per unit leaf area (calculated as : ((mass of Nitrogen/SLA)*1000) in kg m^2)"

# --- 3: insert the new %P row right after the stem N row ---
$ws.Rows(51).Insert()
$ws.Range("A51").Value = "%P"
$ws.Range("B51").Value = "Total leaf phosphorous (%)"

# --- 4: left-align the numeric pi constant ---
$ws.Range("B34").HorizontalAlignment = -4131

# --- 5: widen column B to fit the (now much longer) description text ---
$ws.Columns("B").ColumnWidth = 182

# --- 6: leave the selection where the author's editing session ended ---
$ws.Range("B55").Select()
